# Correct data update following autospc v0.0.0.9030
# Applies corrected statistical summary values to Sheet1 per updated autospc results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 15).Value = 9.478260869565217
$ws.Cells.Item(2, 16).Value = 2.073527648287469
$ws.Cells.Item(2, 18).Value = 33.5
$ws.Cells.Item(2, 19).Value = 53.5
$ws.Cells.Item(2, 20).Value = 65.75
$ws.Cells.Item(2, 21).Value = 135
$ws.Cells.Item(2, 22).Value = 54.93478260869565
$ws.Cells.Item(2, 23).Value = 24.96611229817378
$ws.Cells.Item(2, 25).Value = 97.25
$ws.Cells.Item(2, 26).Value = 128.5
$ws.Cells.Item(2, 28).Value = 223
$ws.Cells.Item(2, 29).Value = 129.2826086956522
$ws.Cells.Item(2, 30).Value = 43.77399686951566
$ws.Cells.Item(2, 32).Value = 7
$ws.Cells.Item(2, 35).Value = 16
$ws.Cells.Item(2, 36).Value = 8.869565217391305
$ws.Cells.Item(2, 37).Value = 3.159526656898494
$ws.Cells.Item(2, 43).Value = 0.5008989865969271
$ws.Cells.Item(2, 44).Value = 0.105439596677146
$ws.Cells.Item(2, 50).Value = 42.21806164849643
$ws.Cells.Item(2, 51).Value = 8.578303307222772
$ws.Cells.Item(5, 14).Value = 17
$ws.Cells.Item(5, 15).Value = 14.39130434782609
$ws.Cells.Item(5, 16).Value = 1.819380369119177
$ws.Cells.Item(5, 20).Value = 53.75
$ws.Cells.Item(5, 21).Value = 98
$ws.Cells.Item(5, 22).Value = 41.78260869565217
$ws.Cells.Item(5, 23).Value = 19.95652279161796
$ws.Cells.Item(5, 25).Value = 81.75
$ws.Cells.Item(5, 28).Value = 190
$ws.Cells.Item(5, 29).Value = 106.1086956521739
$ws.Cells.Item(5, 30).Value = 39.39443178976754
$ws.Cells.Item(5, 36).Value = 8.195652173913043
$ws.Cells.Item(5, 37).Value = 3.074045942094275
$ws.Cells.Item(5, 42).Value = 0.8947368421052632
$ws.Cells.Item(5, 43).Value = 0.7611147433801896
$ws.Cells.Item(5, 44).Value = 0.08751033657951186
$ws.Cells.Item(5, 45).Value = 23.77777777777778
$ws.Cells.Item(5, 50).Value = 28.05624396915957
$ws.Cells.Item(5, 51).Value = 3.784215627775039
$ws.Cells.Item(6, 11).Value = 24
$ws.Cells.Item(6, 13).Value = 30.5
$ws.Cells.Item(6, 15).Value = 27.15217391304348
$ws.Cells.Item(6, 16).Value = 4.50169050372239
$ws.Cells.Item(6, 18).Value = 18.5
$ws.Cells.Item(6, 20).Value = 30
$ws.Cells.Item(6, 22).Value = 25.8695652173913
$ws.Cells.Item(6, 23).Value = 11.24595472090036
$ws.Cells.Item(6, 26).Value = 17
$ws.Cells.Item(6, 27).Value = 23.5
$ws.Cells.Item(6, 29).Value = 16.28260869565218
$ws.Cells.Item(6, 30).Value = 11.61160921660016
$ws.Cells.Item(6, 33).Value = 2
$ws.Cells.Item(6, 36).Value = 1.804347826086957
$ws.Cells.Item(6, 37).Value = 1.25820975494534
$ws.Cells.Item(6, 39).Value = 1.263157894736842
$ws.Cells.Item(6, 43).Value = 1.438460281137627
$ws.Cells.Item(6, 44).Value = 0.2366442513354118
$ws.Cells.Item(6, 46).Value = 13.375
$ws.Cells.Item(6, 50).Value = 15.53341553010698
$ws.Cells.Item(6, 51).Value = 2.557282365008493
$ws.Cells.Item(7, 15).Value = 8.891304347826088
$ws.Cells.Item(7, 16).Value = 2.253178109738087
$ws.Cells.Item(7, 18).Value = 21
$ws.Cells.Item(7, 22).Value = 34.5
$ws.Cells.Item(7, 23).Value = 19.11689189056515
$ws.Cells.Item(7, 27).Value = 75.5
$ws.Cells.Item(7, 29).Value = 54.71739130434783
$ws.Cells.Item(7, 30).Value = 37.82454408536474
$ws.Cells.Item(7, 33).Value = 5
$ws.Cells.Item(7, 36).Value = 4.565217391304348
$ws.Cells.Item(7, 37).Value = 2.613232769434513
$ws.Cells.Item(7, 43).Value = 0.4704151683556718
$ws.Cells.Item(7, 44).Value = 0.1165200850931687
$ws.Cells.Item(7, 50).Value = 45.55266786353743
$ws.Cells.Item(7, 51).Value = 12.01243848599598
$ws.Cells.Item(10, 22).Value = 25.17391304347826
$ws.Cells.Item(10, 23).Value = 10.85317024411882
$ws.Cells.Item(10, 29).Value = 44.89130434782609
$ws.Cells.Item(10, 30).Value = 33.05028375125768
$ws.Cells.Item(11, 14).Value = 32
$ws.Cells.Item(11, 15).Value = 17.3695652173913
$ws.Cells.Item(11, 16).Value = 5.94552565904133
$ws.Cells.Item(11, 22).Value = 33.93478260869565
$ws.Cells.Item(11, 23).Value = 20.80213896471337
$ws.Cells.Item(11, 28).Value = 44
$ws.Cells.Item(11, 29).Value = 13.23913043478261
$ws.Cells.Item(11, 30).Value = 10.89991800704879
$ws.Cells.Item(11, 35).Value = 5
$ws.Cells.Item(11, 36).Value = 1.434782608695652
$ws.Cells.Item(11, 37).Value = 1.186070522786969
$ws.Cells.Item(11, 41).Value = 1.157894736842105
$ws.Cells.Item(11, 42).Value = 1.684210526315789
$ws.Cells.Item(11, 43).Value = 0.9211343576332135
$ws.Cells.Item(11, 44).Value = 0.3160126604494409
$ws.Cells.Item(11, 45).Value = 12.96969696969697
$ws.Cells.Item(11, 46).Value = 18.60869565217391
$ws.Cells.Item(11, 50).Value = 26.19468100350008
$ws.Cells.Item(11, 51).Value = 10.27406184135486
$ws.Cells.Item(12, 15).Value = 1.077777777777778
$ws.Cells.Item(12, 16).Value = 0.7451620084486273
$ws.Cells.Item(12, 18).Value = 14
$ws.Cells.Item(12, 19).Value = 21
$ws.Cells.Item(12, 20).Value = 28.75
$ws.Cells.Item(12, 22).Value = 21.84814814814815
$ws.Cells.Item(12, 23).Value = 11.36759466438732
$ws.Cells.Item(12, 25).Value = 23.25
$ws.Cells.Item(12, 26).Value = 36
$ws.Cells.Item(12, 27).Value = 46
$ws.Cells.Item(12, 29).Value = 35.52592592592593
$ws.Cells.Item(12, 30).Value = 16.09443766835809
$ws.Cells.Item(12, 37).Value = 1.070012307863037
$ws.Cells.Item(12, 43).Value = 0.5358024691358024
$ws.Cells.Item(12, 44).Value = 0.3738593751361511
$ws.Cells.Item(12, 45).Value = 24.25
$ws.Cells.Item(12, 47).Value = 47.5
$ws.Cells.Item(12, 50).Value = 47.90246913580247
$ws.Cells.Item(12, 51).Value = 19.15548415972048
$ws.Cells.Item(13, 22).Value = 31.90370370370371
$ws.Cells.Item(13, 23).Value = 21.68269797410013
$ws.Cells.Item(13, 29).Value = 59.24814814814815
$ws.Cells.Item(13, 30).Value = 21.36515368765059
$ws.Cells.Item(13, 33).Value = 2
$ws.Cells.Item(13, 36).Value = 2.662962962962963
$ws.Cells.Item(13, 37).Value = 1.244063418980206
$ws.Cells.Item(15, 15).Value = 1.67037037037037
$ws.Cells.Item(15, 16).Value = 0.5708630489707315
$ws.Cells.Item(15, 20).Value = 26
$ws.Cells.Item(15, 22).Value = 18.54444444444444
$ws.Cells.Item(15, 23).Value = 11.32831545742129
$ws.Cells.Item(15, 27).Value = 48
$ws.Cells.Item(15, 28).Value = 80
$ws.Cells.Item(15, 29).Value = 36.06666666666667
$ws.Cells.Item(15, 30).Value = 16.26486166270907
$ws.Cells.Item(15, 36).Value = 2.222222222222222
$ws.Cells.Item(15, 37).Value = 1.088409136572718
$ws.Cells.Item(15, 43).Value = 0.8493827160493828
$ws.Cells.Item(15, 44).Value = 0.2672191257420775
$ws.Cells.Item(15, 45).Value = 24.25
$ws.Cells.Item(15, 50).Value = 33.88302469135802
$ws.Cells.Item(15, 51).Value = 6.463459304643556
$ws.Cells.Item(16, 15).Value = 3.785185185185185
$ws.Cells.Item(16, 16).Value = 1.323986237336714
$ws.Cells.Item(16, 19).Value = 6
$ws.Cells.Item(16, 20).Value = 10
$ws.Cells.Item(16, 22).Value = 7.159259259259259
$ws.Cells.Item(16, 23).Value = 6.608587364827454
$ws.Cells.Item(16, 27).Value = 16
$ws.Cells.Item(16, 28).Value = 36
$ws.Cells.Item(16, 29).Value = 7.892592592592592
$ws.Cells.Item(16, 30).Value = 8.984473427541825
$ws.Cells.Item(16, 36).Value = 0.774074074074074
$ws.Cells.Item(16, 37).Value = 0.8297859233004039
$ws.Cells.Item(16, 42).Value = 4
$ws.Cells.Item(16, 43).Value = 1.957407407407407
$ws.Cells.Item(16, 44).Value = 0.6931963174436352
$ws.Cells.Item(16, 45).Value = 11.875
$ws.Cells.Item(16, 48).Value = 23.75
$ws.Cells.Item(16, 50).Value = 19.64593915343915
$ws.Cells.Item(16, 51).Value = 6.079455385210363
$ws.Cells.Item(17, 15).Value = 1.287179487179487
$ws.Cells.Item(17, 16).Value = 0.6731903911003474
$ws.Cells.Item(17, 18).Value = 14.5
$ws.Cells.Item(17, 22).Value = 23.81538461538462
$ws.Cells.Item(17, 23).Value = 11.5587593076745
$ws.Cells.Item(17, 26).Value = 36
$ws.Cells.Item(17, 29).Value = 36.31794871794872
$ws.Cells.Item(17, 30).Value = 14.1170159771975
$ws.Cells.Item(17, 36).Value = 2.435897435897436
$ws.Cells.Item(17, 37).Value = 0.9684420423682468
$ws.Cells.Item(17, 43).Value = 0.6025641025641025
$ws.Cells.Item(17, 44).Value = 0.3197744601253966
$ws.Cells.Item(17, 50).Value = 44.18504273504274
$ws.Cells.Item(17, 51).Value = 16.10794095059001
$ws.Cells.Item(18, 19).Value = 34
$ws.Cells.Item(18, 22).Value = 30.45714285714286
$ws.Cells.Item(18, 23).Value = 26.28282114134587
$ws.Cells.Item(18, 29).Value = 46.90357142857143
$ws.Cells.Item(18, 30).Value = 35.00764875095368
$ws.Cells.Item(18, 33).Value = 2.5
$ws.Cells.Item(18, 36).Value = 2.242857142857143
$ws.Cells.Item(18, 37).Value = 1.773445822980102
$ws.Cells.Item(20, 15).Value = 1.641025641025641
$ws.Cells.Item(20, 16).Value = 0.5871115512114823
$ws.Cells.Item(20, 18).Value = 15
$ws.Cells.Item(20, 19).Value = 24
$ws.Cells.Item(20, 22).Value = 23.03589743589744
$ws.Cells.Item(20, 23).Value = 10.70425985103587
$ws.Cells.Item(20, 25).Value = 25
$ws.Cells.Item(20, 27).Value = 46
$ws.Cells.Item(20, 29).Value = 35.67692307692307
$ws.Cells.Item(20, 30).Value = 15.21743542779817
$ws.Cells.Item(20, 36).Value = 2.384615384615385
$ws.Cells.Item(20, 37).Value = 1.084458737513946
$ws.Cells.Item(20, 39).Value = 0.5
$ws.Cells.Item(20, 43).Value = 0.7837606837606838
$ws.Cells.Item(20, 44).Value = 0.2893560210584588
$ws.Cells.Item(20, 48).Value = 37.25
$ws.Cells.Item(20, 50).Value = 36.77435897435898
$ws.Cells.Item(20, 51).Value = 11.01836489330516
$ws.Cells.Item(21, 15).Value = 4.492307692307692
$ws.Cells.Item(21, 16).Value = 1.503841235482809
$ws.Cells.Item(21, 22).Value = 10.64615384615385
$ws.Cells.Item(21, 23).Value = 7.508271938467769
$ws.Cells.Item(21, 27).Value = 17
$ws.Cells.Item(21, 29).Value = 9.805128205128206
$ws.Cells.Item(21, 30).Value = 8.817331477873292
$ws.Cells.Item(21, 34).Value = 2
$ws.Cells.Item(21, 36).Value = 1.015384615384615
$ws.Cells.Item(21, 37).Value = 0.8994667554123335
$ws.Cells.Item(21, 41).Value = 2.5
$ws.Cells.Item(21, 43).Value = 2.152136752136752
$ws.Cells.Item(21, 44).Value = 0.7309550184070249
$ws.Cells.Item(21, 50).Value = 18.47757224257224
$ws.Cells.Item(21, 51).Value = 9.143566879842918
